# Refresh the cryptocurrency price/volume snapshot (rows 2-51, columns D "Price"
# and E "Volume(1h)") to the latest scraped values. Numeric-looking price strings
# are written with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cells) instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.549.08'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.077.50'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'523.22"
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").Value = "'140.32"
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.075.07'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = "'0.440"
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").Value = "'7.16"
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("D13").Value = '3.606.57'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").Value = "'25.34"
$ws.Range("E15").Value = '  -5.50%  '
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '57.553.39'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '3.073.06'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = "'6.06"
$ws.Range("E19").Value = '  -1.68%  '
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").Value = "'7.92"
$ws.Range("E21").Value = '  -2.51%  '
$ws.Range("D22").Value = "'338.70"
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("D24").Value = "'0.508"
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = "'66.89"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("E26").Value = '  -1.80%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '0.0₃0909'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = "'6.35"
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("D31").Value = "'7.20"
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  +2.98%  '
$ws.Range("D33").Value = "'20.83"
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("D35").Value = "'158.46"
$ws.Range("E35").Value = '  +1.95%  '
$ws.Range("D36").Value = "'4.58"
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("E37").Value = '  +1.40%  '
$ws.Range("E38").Value = '  -5.96%  '
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("D40").Value = "'0.0662"
$ws.Range("E40").Value = '  -2.04%  '
$ws.Range("D41").Value = "'1.57"
$ws.Range("E41").Value = '  +12.90%  '
$ws.Range("D42").Value = "'3.97"
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("D43").Value = "'0.681"
$ws.Range("E43").Value = '  +3.75%  '
$ws.Range("D44").Value = '3.115.98'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = "'36.80"
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '2.273.96'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("D49").Value = "'0.988"
$ws.Range("E49").Value = '  +4.42%  '
$ws.Range("D50").Value = "'6.06"
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").Value = "'20.44"
$ws.Range("E51").Value = '  -0.51%  '
